# Aula 03.1 - Arduino Tinkercad: fix SDA pin note from A6 to A4
# (commit: "arduino ajuste 10/09/2024 2024.2")
#
# Slide 10, shape "CaixaDeTexto 1" has a bullet line built from runs:
#   "SDA => "  (plain)  +  "A6" (blue, 0070C0)  +  " " (plain)
# It must become:
#   "SDA "  +  "=> "  +  "A4" (blue, 0070C0)  +  " "
# i.e. the pin reference changes from A6 to A4, with the "SDA => " run
# split in two (formatting identical on both halves).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item("CaixaDeTexto 1")
$tr = $shape.TextFrame.TextRange

$full = $tr.Text
$idx = $full.IndexOf("SDA => A6")

if ($idx -ge 0) {
    # Split "SDA => " (7 chars) into "SDA " (4 chars) + "=> " (3 chars).
    # Re-assigning the same text on a sub-range forces PowerPoint to turn
    # it into its own run while keeping the original character formatting.
    $sdaPart = $tr.Characters($idx + 1, 4)
    $sdaPart.Text = "SDA "

    # Replace the pin number "A6" -> "A4" (keeps the blue-colored run).
    $pinPart = $tr.Characters($idx + 8, 2)
    $pinPart.Text = "A4"
}
